# Workbook: faturamento_diario.xlsx
# Commit: "atualizei dados da bibi e add"
#
# 1) A new daily-sales row is inserted into the table (sorted by month
#    desc, then day asc): day 10 of 07/2025, total_venda = 12701.99.
#    It is inserted right after day 9 of 07/2025 (i.e. immediately before
#    the first row of the 06/2025 block), pushing every row below it down
#    by one. The sheet's used range grows from A1:E70 to A1:E71.
#
# 2) The existing 06/2025 "day 30" total_venda figure is corrected from
#    111900.66 to 109507.06 (same row, now shifted down to row 30 by the
#    insertion above).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new 07/2025 row (day 10) at row 9 -------------------------
# Before insert, row 9 holds day 2 of 06/2025; inserting pushes it (and
# everything after) down by one row.
$ws.Rows("9:9").Insert()

$ws.Range("A9").Value = 10
$ws.Range("B9").Value = 12701.99
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = 2025
$ws.Range("E9").Value = "07/2025"

# --- Correct the 06/2025 day-30 total ------------------------------------
# That row has shifted from row 29 to row 30 because of the insertion.
$ws.Range("B30").Value = 109507.06
